# Generate Report for Handback
# Updates the localization-status report after a handback run:
#  - Overview row for cceb51ea... now reports a failed handback transform
#    instead of "Ready for handoff".
#  - The zh-cn and de-de detail sheets get an Error Detail message in
#    column P (row 3 = cceb51ea... file), and that column is widened so
#    the message is readable.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = "Handback transform failed"
$ws1.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667
$wsZhCn.Range("P3").Value = "Handback file name: z5ercxj4.gpc is different with handoff file name: cceb51ea-1423-4eed-b9b8-0fa25ab789f1.5ffb3d71a8dd62e6266fd3476fc389d11f1c9764.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
$wsDeDe.Range("P3").Value = "Handback file name: z5ercxj4.gpc is different with handoff file name: cceb51ea-1423-4eed-b9b8-0fa25ab789f1.5ffb3d71a8dd62e6266fd3476fc389d11f1c9764.de-de."
